# mise en place survey
# Refresh market cap figures, and re-sort a couple of rows whose
# relative ranking changed (OKB now above Cosmos, Kaspa now above Monero).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Market Cap (column C) values that kept their row/rank ---
$ws.Range("C2").Value  = 734020108300.9487
$ws.Range("C3").Value  = 246066347323.4598
$ws.Range("C4").Value  = 38992486587.04223
$ws.Range("C5").Value  = 33307298629.42934
$ws.Range("C6").Value  = 24253458545.67005
$ws.Range("C7").Value  = 13677197642.06798
$ws.Range("C8").Value  = 11359429835.55813
$ws.Range("C9").Value  = 9119366109.749874
$ws.Range("C10").Value = 8257079728.487464
$ws.Range("C11").Value = 8093285119.901435
$ws.Range("C12").Value = 7637931566.491719
$ws.Range("C13").Value = 7592349038.162902
$ws.Range("C14").Value = 6727678885.928888
$ws.Range("C15").Value = 6131783021.903178
$ws.Range("C16").Value = 5185208942.700999
$ws.Range("C17").Value = 5060123103.839664
$ws.Range("C18").Value = 4467550118.133206
$ws.Range("C19").Value = 3651225093.29196

# --- Rows 20/21 swap: OKB now ranks above Cosmos ---
$ws.Range("A20").Value = "OKB"
$ws.Range("B20").Value = "OKB-USD"
$ws.Range("C20").Value = 3454728662.767564

$ws.Range("A21").Value = "Cosmos"
$ws.Range("B21").Value = "ATOM-USD"
$ws.Range("C21").Value = 3453696439.300338

$ws.Range("C22").Value = 3366126614.056954

# --- Rows 23/24 swap: Kaspa now ranks above Monero ---
$ws.Range("A23").Value = "Kaspa"
$ws.Range("B23").Value = "KAS-USD"
$ws.Range("C23").Value = 3015161535.871074

$ws.Range("A24").Value = "Monero"
$ws.Range("B24").Value = "XMR-USD"
$ws.Range("C24").Value = 3004347745.952567

$ws.Range("C25").Value = 2841200082.369965
$ws.Range("C26").Value = 2440675830.95558
